# Survey results table: rename the ranking question header, grow the table
# from 3 submissions to 10, and (re)populate every data row with the final
# set of responses.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Rename the last table column header (ranking question wording changed)
$ws.Range("F1").Value = "Rank each of the following pitches, starting with your most-preferred project at the top of the list."

# Grow the table from 3 data rows to 10 data rows (A1:F4 -> A1:F11)
$lo.Resize($ws.Range("A1:F11"))

# Carry the existing date/time cell formatting down into the newly added rows
# before writing their values, so the new Start/Completion time cells render
# the same way as the existing ones.
$ws.Range("B2:C2").Copy()
$ws.Range("B5:C11").PasteSpecial(-4122)

# Populate every data row (rows 2-4 are replaced, rows 5-11 are newly added)
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 44023.7712962963
$ws.Range("C2").Value = 44023.773125
$ws.Range("D2").Value = "ssmeltze@uwo.ca"
$ws.Range("E2").Value = "Sandra Christine Smeltzer"
$ws.Range("F2").Value = "De Groot;Lee;Smye;Beveridge;McNair;Tang;Esses;Bitsuamlak;Hill;Petrella;"

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 44024.4625925926
$ws.Range("C3").Value = 44024.4737384259
$ws.Range("D3").Value = "abottere@uwo.ca"
$ws.Range("E3").Value = "Andrew Botterell"
$ws.Range("F3").Value = "Lee;McNair;Smye;De Groot;Tang;Beveridge;Hill;Esses;Bitsuamlak;Petrella;"

$ws.Range("A4").Value = 6
$ws.Range("B4").Value = 44025.506087963
$ws.Range("C4").Value = 44025.5067708333
$ws.Range("D4").Value = "litchfi@uwo.ca"
$ws.Range("E4").Value = "David William Litchfield"
$ws.Range("F4").Value = "Tang;McNair;Smye;Lee;Petrella;De Groot;Beveridge;Hill;Esses;Bitsuamlak;"

$ws.Range("A5").Value = 7
$ws.Range("B5").Value = 44025.6872453704
$ws.Range("C5").Value = 44025.6879976852
$ws.Range("D5").Value = "joramcar@uwo.ca"
$ws.Range("E5").Value = "Janis Cardy"
$ws.Range("F5").Value = "Tang;Smye;Esses;De Groot;Lee;Hill;McNair;Petrella;Beveridge;Bitsuamlak;"

$ws.Range("A6").Value = 8
$ws.Range("B6").Value = 44025.6959837963
$ws.Range("C6").Value = 44025.6978703704
$ws.Range("D6").Value = "jburkell@uwo.ca"
$ws.Range("E6").Value = "Jacquelyn Burkell"
$ws.Range("F6").Value = "Smye;Tang;Lee;Esses;McNair;De Groot;Hill;Bitsuamlak;Beveridge;Petrella;"

$ws.Range("A7").Value = 9
$ws.Range("B7").Value = 44025.7209953704
$ws.Range("C7").Value = 44025.7305208333
$ws.Range("D7").Value = "mcapretz@uwo.ca"
$ws.Range("E7").Value = "Miriam Capretz"
$ws.Range("F7").Value = "Smye;Tang;De Groot;Petrella;Bitsuamlak;Esses;McNair;Beveridge;Hill;Lee;"

$ws.Range("A8").Value = 10
$ws.Range("B8").Value = 44026.3932060185
$ws.Range("C8").Value = 44026.4142476852
$ws.Range("D8").Value = "ascully2@uwo.ca"
$ws.Range("E8").Value = "Abbey Baran"
$ws.Range("F8").Value = "Petrella;Lee;McNair;De Groot;Smye;Beveridge;Bitsuamlak;Tang;Hill;Esses;"

$ws.Range("A9").Value = 11
$ws.Range("B9").Value = 44026.4249768518
$ws.Range("C9").Value = 44026.4394675926
$ws.Range("D9").Value = "bneff@uwo.ca"
$ws.Range("E9").Value = "Bryan Neff"
$ws.Range("F9").Value = "Bitsuamlak;Lee;De Groot;Hill;McNair;Tang;Beveridge;Smye;Esses;Petrella;"

$ws.Range("A10").Value = 12
$ws.Range("B10").Value = 44026.4753009259
$ws.Range("C10").Value = 44026.5234259259
$ws.Range("D10").Value = "eabrams3@uwo.ca"
$ws.Range("E10").Value = "Emily Ansari"
$ws.Range("F10").Value = "McNair;Esses;De Groot;Lee;Bitsuamlak;Tang;Smye;Petrella;Hill;Beveridge;"

$ws.Range("A11").Value = 13
$ws.Range("B11").Value = 44026.6851967593
$ws.Range("C11").Value = 44026.8601736111
$ws.Range("D11").Value = "kenm@uwo.ca"
$ws.Range("E11").Value = "Ken McRae"
$ws.Range("F11").Value = "Lee;Tang;Hill;De Groot;Esses;Smye;McNair;Beveridge;Petrella;Bitsuamlak;"
